$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 243.06229999999999
$ws.Range("C11").Value = 251

$ws.Range("A11:B11").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()

$ws.Range("A14:D15").Select()
